$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.212.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.415.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.844.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.139.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.472.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "328.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0772"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("E32").Value = "  +7.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.402"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("E35").Value = "  +4.06%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "324.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.81%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "146.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0971"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0517"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.577"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("E51").Value = "  -0.80%  "
